# Replace every cell whose value is exactly "X" with "-" across all worksheets.
$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count
$totalChanged = 0

for ($s = 1; $s -le $sheetCount; $s++) {
  $ws = $wb.Worksheets.Item($s)
  $used = $ws.UsedRange
  $rowCount = $used.Rows.Count
  $colCount = $used.Columns.Count

  for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
      $cell = $used.Cells.Item($r, $c)
      if ($cell.Value2 -eq "X") {
        $cell.Value = "-"
        $totalChanged = $totalChanged + 1
      }
    }
  }
}

"Total cells changed: " + $totalChanged
